$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto market data.
# Some Price values are plain decimal-looking strings (e.g. "586.04") that the Excel
# COM layer would otherwise auto-convert to numbers; force those cells to Text format
# first so the values are stored as literal text, matching the source data feed.

$ws.Range('D2').Value = '70.442.05'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '3.576.17'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '586.04'
$ws.Range('E5').Value = '  +1.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.86'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').Value = '3.564.39'
$ws.Range('E7').Value = '  -0.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.619'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.207'
$ws.Range('E10').Value = '  +13.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.649'
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.19'
$ws.Range('E12').Value = '  -1.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000314'
$ws.Range('E13').Value = '  +3.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.54'
$ws.Range('E14').Value = '  -0.11%  '
$ws.Range('D15').Value = '4.145.68'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.59'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('D17').Value = '70.452.53'
$ws.Range('E17').Value = '  +0.55%  '
$ws.Range('D18').Value = '3.570.56'
$ws.Range('E18').Value = '  -0.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.41'
$ws.Range('E19').Value = '  -1.70%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '566.15'
$ws.Range('E20').Value = '  +15.85%  '
$ws.Range('E21').Value = '  -0.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.01'
$ws.Range('E22').Value = '  -1.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '17.80'
$ws.Range('E23').Value = '  -6.47%  '
$ws.Range('E24').Value = '  +5.71%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.91'
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '95.30'
$ws.Range('E26').Value = '  +0.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.59'
$ws.Range('E27').Value = '  -1.69%  '
$ws.Range('E28').Value = '  -0.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.14'
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.19'
$ws.Range('E30').Value = '  +1.84%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.33'
$ws.Range('E31').Value = '  -6.12%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.47'
$ws.Range('E32').Value = '  +2.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.86'
$ws.Range('E33').Value = '  -2.18%  '
$ws.Range('E34').Value = '  -0.45%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.36'
$ws.Range('E35').Value = '  +2.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '565.57'
$ws.Range('E36').Value = '  -1.55%  '
$ws.Range('E37').Value = '  -0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '37.83'
$ws.Range('E38').Value = '  -2.50%  '
$ws.Range('E39').Value = '  +0.08%  '
$ws.Range('D40').Value = '0.0₃0773'
$ws.Range('E40').Value = '  -2.18%  '
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('D42').Value = '3.367.84'
$ws.Range('E42').Value = '  +4.78%  '
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('E44').Value = '  -3.28%  '
$ws.Range('E45').Value = '  +2.15%  '
$ws.Range('E46').Value = '  +0.21%  '
$ws.Range('E47').Value = '  -4.08%  '
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('E49').Value = '  +0.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.998'
$ws.Range('E50').Value = '  -0.15%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.48'
$ws.Range('E51').Value = '  -10.74%  '
